$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.431.90'
$ws.Range("E2").Value = '  -0.72%  '

$ws.Range("D3").Value = '3.125.89'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").Value = '''215.60'
$ws.Range("E5").Value = '  -1.71%  '

$ws.Range("D6").Value = '''620.75'
$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").Value = '''1.12'
$ws.Range("E7").Value = '  +25.80%  '

$ws.Range("D8").Value = '''0.362'
$ws.Range("E8").Value = '  -4.27%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").Value = '3.124.39'
$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("D11").Value = '''0.735'
$ws.Range("E11").Value = '  +5.61%  '

$ws.Range("E12").Value = '  +5.67%  '

$ws.Range("D13").Value = '''0.0000245'
$ws.Range("E13").Value = '  -3.74%  '

$ws.Range("D14").Value = '''5.63'
$ws.Range("E14").Value = '  +4.50%  '

$ws.Range("D15").Value = '''35.09'
$ws.Range("E15").Value = '  +6.23%  '

$ws.Range("D16").Value = '90.247.25'
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").Value = '3.712.22'
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").Value = '3.138.02'

$ws.Range("E19").Value = '  +3.95%  '

$ws.Range("D20").Value = '''14.47'
$ws.Range("E20").Value = '  +4.69%  '

$ws.Range("D21").Value = '''0.0000211'
$ws.Range("E21").Value = '  -10.93%  '

$ws.Range("D22").Value = '''461.33'
$ws.Range("E22").Value = '  +7.04%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '''9.03'
$ws.Range("E23").Value = '  +5.34%  '

$ws.Range("B24").Value = 'Polkadot'
$ws.Range("C24").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D24").Value = '''5.36'
$ws.Range("E24").Value = '  +4.10%  '

$ws.Range("D25").Value = '''94.75'
$ws.Range("E25").Value = '  +13.32%  '

$ws.Range("D26").Value = '''5.76'
$ws.Range("E26").Value = '  +3.13%  '

$ws.Range("D27").Value = '''12.23'
$ws.Range("E27").Value = '  +2.78%  '

$ws.Range("D28").Value = '3.306.15'
$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("D30").Value = '''0.163'
$ws.Range("E30").Value = '  -2.88%  '

$ws.Range("D31").Value = '''0.219'
$ws.Range("E31").Value = '  +50.44%  '

$ws.Range("D32").Value = '''9.19'
$ws.Range("E32").Value = '  +5.73%  '

$ws.Range("D33").Value = '''26.55'
$ws.Range("E33").Value = '  +15.45%  '

$ws.Range("D34").Value = '''515.79'
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("D35").Value = '''0.145'
$ws.Range("E35").Value = '  +3.74%  '

$ws.Range("E36").Value = '  +4.86%  '

$ws.Range("D37").Value = '''6.99'
$ws.Range("E37").Value = '  +0.69%  '

$ws.Range("D38").Value = '''1.32'
$ws.Range("E38").Value = '  +2.54%  '

$ws.Range("D39").Value = '''3.58'
$ws.Range("E39").Value = '  -9.19%  '

$ws.Range("D40").Value = '''0.0922'
$ws.Range("E40").Value = '  +28.43%  '

$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '''22.22'
$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").Value = '''0.426'
$ws.Range("E42").Value = '  +14.49%  '

$ws.Range("E43").Value = '  -24.78%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("E45").Value = '  +5.73%  '

$ws.Range("D47").Value = '''0.720'
$ws.Range("E47").Value = '  +18.57%  '

$ws.Range("D48").Value = '''4.68'
$ws.Range("E48").Value = '  +11.27%  '

$ws.Range("D49").Value = '''150.55'
$ws.Range("E49").Value = '  +6.44%  '

$ws.Range("D50").Value = '''1.36'
$ws.Range("E50").Value = '  +7.93%  '

$ws.Range("E51").Value = '  +3.29%  '
